$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'55.098.61"
$ws.Range("E2").Value = "  -1.91%  "
$ws.Range("D3").Value = "'2.346.00"
$ws.Range("E3").Value = "  -5.04%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'475.29"
$ws.Range("E5").Value = "  -2.39%  "
$ws.Range("D6").Value = "'145.51"
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.603"
$ws.Range("E8").Value = "  +18.60%  "
$ws.Range("D9").Value = "'2.350.92"
$ws.Range("E9").Value = "  -5.20%  "
$ws.Range("D10").Value = "'0.0967"
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("D11").Value = "'5.43"
$ws.Range("E11").Value = "  -5.96%  "
$ws.Range("D12").Value = "'0.324"
$ws.Range("E12").Value = "  -2.09%  "
$ws.Range("D13").Value = "'0.125"
$ws.Range("E13").Value = "  +1.21%  "
$ws.Range("D14").Value = "'2.751.23"
$ws.Range("E14").Value = "  -5.51%  "
$ws.Range("D15").Value = "'55.094.83"
$ws.Range("E15").Value = "  -2.01%  "
$ws.Range("D16").Value = "'20.03"
$ws.Range("E16").Value = "  -4.60%  "
$ws.Range("D17").Value = "'0.0000130"
$ws.Range("E17").Value = "  -3.61%  "
$ws.Range("D18").Value = "'2.349.58"
$ws.Range("E18").Value = "  -5.29%  "
$ws.Range("D19").Value = "'4.58"
$ws.Range("E19").Value = "  +1.96%  "
$ws.Range("D20").Value = "'316.00"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").Value = "'9.61"
$ws.Range("E21").Value = "  -3.93%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "'5.62"
$ws.Range("E23").Value = "  -2.73%  "
$ws.Range("D24").Value = "'56.54"
$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "'0.394"
$ws.Range("E26").Value = "  -4.06%  "
$ws.Range("D27").Value = "'0.154"
$ws.Range("E27").Value = "  -5.30%  "
$ws.Range("D28").Value = "'2.440.93"
$ws.Range("E28").Value = "  -5.50%  "
$ws.Range("D29").Value = "'7.09"
$ws.Range("E29").Value = "  -6.74%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").Value = "0.0₃0749"
$ws.Range("E31").Value = "  -4.79%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'18.20"
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "'145.03"
$ws.Range("E33").Value = "  -2.71%  "
$ws.Range("D34").Value = "'1.47"
$ws.Range("E34").Value = "  -2.16%  "
$ws.Range("D35").Value = "'5.10"
$ws.Range("E35").Value = "  -1.47%  "
$ws.Range("D36").Value = "'3.59"
$ws.Range("E36").Value = "  -3.43%  "
$ws.Range("D37").Value = "'1.09"
$ws.Range("E37").Value = "  -4.08%  "
$ws.Range("D38").Value = "'0.816"
$ws.Range("E38").Value = "  -4.72%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "'33.72"
$ws.Range("E39").Value = "  -1.30%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "'0.0992"
$ws.Range("E40").Value = "  +7.25%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("D42").Value = "'1.34"
$ws.Range("E42").Value = "  +1.14%  "
$ws.Range("D43").Value = "'3.40"
$ws.Range("E43").Value = "  -3.39%  "
$ws.Range("E44").Value = "  -4.69%  "
$ws.Range("D45").Value = "'0.0520"
$ws.Range("E45").Value = "  -6.28%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "'250.84"
$ws.Range("E47").Value = "  -2.45%  "
$ws.Range("D48").Value = "'0.0221"
$ws.Range("E48").Value = "  -2.66%  "
$ws.Range("D49").Value = "'4.36"
$ws.Range("E49").Value = "  -8.98%  "
$ws.Range("D50").Value = "'16.75"
$ws.Range("E50").Value = "  -4.39%  "
$ws.Range("D51").Value = "'1.771.64"
$ws.Range("E51").Value = "  -4.62%  "

# Reset number format/style on forced-text cells so they match plain default styling
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
